# Fix missing Spanish accents / naming in the "comida" sheet's text columns,
# then re-sync the AutoFilter range / _FilterDatabase defined name with the
# sheet's actual data extent (A1:I132), and finally move the selection to K5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

# --- Spelling / accent corrections -----------------------------------------
# (substring replace covers compound phrases like "atun c/ arroz",
#  "atun c/ ensalada", "tortellini de salmon al curry" and
#  "bife c/ pure calabaza" automatically)
$cells.Replace("salmon", "salmón")
$cells.Replace("atun", "atún")
$cells.Replace("morron", "morrón")
$cells.Replace("jamon", "jamón")
$cells.Replace("marron", "marrón")
$cells.Replace("lacteo", "lácteo")
$cells.Replace("burger", "hamburguesa")
$cells.Replace("pure", "puré")
$cells.Replace("anana", "ananá")

# --- Resync AutoFilter range to the full used data range (A1:I132) ---------
$ws.AutoFilterMode = $false
$ws.Range("A1:I132").AutoFilter()

# --- Resync the hidden _FilterDatabase defined name to match ---------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "comida!_FilterDatabase") {
        $n.RefersTo = "=comida!`$A`$1:`$I`$132"
    }
}

# --- Move the active selection to K5 ----------------------------------------
$ws.Range("K5").Select()
